$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 226, shifting the existing rows 226-319 down to 227-320.
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with the new record's data.
$ws.Range("A226").Value = 10
$ws.Range("B226").Value = "Vega Modelo de Temuco"
$ws.Range("C226").Value = "La Araucanía"
$ws.Range("D226").Value = 44609
$ws.Range("E226").Value = 9
$ws.Range("F226").Value = 100112037
$ws.Range("G226").Value = "Cebollín"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 50
$ws.Range("K226").Value = 8000
$ws.Range("L226").Value = 8000
$ws.Range("M226").Value = 8000
$ws.Range("N226").Value = "$/docena de paquetes"
$ws.Range("O226").Value = "Provincia de Cautín"
$ws.Range("P226").Value = 667
$ws.Range("Q226").Value = 12
$ws.Range("R226").Value = "Hortaliza"
